# Upload validation and error handling (#6327)
#
# - Rename the "RO & CO Hearing Allocation" sheet to "RO Allocations"
# - Drop the "Central Office" allocation row (the sheet no longer tracks
#   Central Office hearings, only Regional Office video hearings), which
#   shifts every row below it up by one
# - Update the sheet's title cell to match the new, RO-only scope

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RO & CO Hearing Allocation")

# Remove the "Central Office" row (row 4: label in column B, count in
# column D) - deleting the whole row shifts all subsequent rows (and their
# formatting, including the special bottom-border on the final row) up by
# one, which is exactly what happened in the source workbook.
$ws.Rows.Item(4).Delete()

# Update the big banner/title cell to drop the "and Central Office
# Hearings" portion of the text now that the row is gone.
$ws.Cells.Item(1, 1).Value = "Allocation of Regional Office Video Hearings"

# Rename the sheet itself to reflect the RO-only scope.
$ws.Name = "RO Allocations"
